$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 15 de Septiembre de 2020 a las 05:13"

# Venezuela overtakes Barein in total cases -> rows 53/54 swap country + stats
$ws.Range("A53").Value = "Venezuela"
$ws.Range("A54").Value = "Barein"

# Refresh numeric country statistics (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes)
$ws.Range("B30").Value = 127619
$ws.Range("C30").Value = 828
$ws.Range("D30").Value = 85198
$ws.Range("E30").Value = 35027
$ws.Range("G30").Value = 50
$ws.Range("H30").Value = 7394

$ws.Range("B33").Value = 106920
$ws.Range("C33").Value = 65
$ws.Range("E33").Value = 4659

$ws.Range("B39").Value = 94306
$ws.Range("C39").Value = 851
$ws.Range("D39").Value = 18737
$ws.Range("E39").Value = 65642
$ws.Range("G39").Value = 2
$ws.Range("H39").Value = 9927

$ws.Range("B50").Value = 68620
$ws.Range("C50").Value = 831
$ws.Range("D50").Value = 18487
$ws.Range("E50").Value = 48046
$ws.Range("G50").Value = 8
$ws.Range("H50").Value = 2087

$ws.Range("B53").Value = 61569
$ws.Range("D53").Value = 49371
$ws.Range("E53").Value = 11704
$ws.Range("H53").Value = 494

$ws.Range("B54").Value = 60965
$ws.Range("D54").Value = 54204
$ws.Range("E54").Value = 6548
$ws.Range("H54").Value = 213

$ws.Range("B74").Value = 28367
$ws.Range("D74").Value = 14814
$ws.Range("E74").Value = 13014
$ws.Range("H74").Value = 539

$ws.Range("B140").Value = 3008
$ws.Range("C140").Value = 34
$ws.Range("D140").Value = 1391
$ws.Range("E140").Value = 1549
$ws.Range("G140").Value = 1
$ws.Range("H140").Value = 68

$ws.Range("B155").Value = 1801
$ws.Range("C155").Value = 3
$ws.Range("D155").Value = 1694
$ws.Range("E155").Value = 83

$ws.Range("B159").Value = 1501
$ws.Range("C159").Value = 21
$ws.Range("D159").Value = 540
$ws.Range("E159").Value = 942

$ws.Range("B172").Value = 648
$ws.Range("C172").Value = 2
$ws.Range("E172").Value = 111

$ws.Range("D185").Value = 301
$ws.Range("E185").Value = 10
